# Tripadvisor New Orleans shard 163 edit
#
# The workbook originally has:
#   Worksheet #1 (rId1): name = "hotel_info"  -> header (9 cols) + 1 data row
#   Worksheet #2 (rId2): name = "review_info" -> header (25 cols), no data rows
#
# After the edit:
#   Worksheet #1 (rId1): name = "review_info" -> header (25 cols) only
#   Worksheet #2 (rId2): name = "hotel_info"  -> header (10 cols, new "State"
#                         column inserted after "Hotel_Name") + 1 data row
#
# i.e. the two sheets swap identity (name + contents), with a new "State"
# column added to the hotel_info table.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- Step 1: swap the sheet names (use a temp name to avoid a collision) ---
$ws1.Name = "__tmp_swap_name__"
$ws2.Name = "hotel_info"
$ws1.Name = "review_info"

# --- Step 2: rebuild worksheet #1 ("review_info") as a header-only sheet ---
$ws1.Cells.Clear()

$reviewHeaders = @(
    "STR","reviewer_ID","reviewer_name","Review_ID","Date_of_scraping","ReviewURL",
    "Tripadvisor_gcode","Tripadvisor_dcode","Tripadvisor_rcode","review_date","review_title",
    "review_content","review_rating","trip_month","trip_purpose","value","rooms","Location",
    "Cleanliness","Sleep Quality","Service","Picture(yes=1)","respondent","response_date",
    "response_text"
)
for ($i = 0; $i -lt $reviewHeaders.Length; $i++) {
    $ws1.Cells.Item(1, $i + 1).Value = $reviewHeaders[$i]
}

# --- Step 3: rebuild worksheet #2 ("hotel_info") with the new State column ---
$ws2.Cells.Clear()

$hotelHeaders = @(
    "STR","Hotel_Name","State","City","Zip","TA_ReviewURL","Tripadvisor_Hotel_Name",
    "English_Reviews_num","Local_Rank","Total_Reviews_num"
)
for ($i = 0; $i -lt $hotelHeaders.Length; $i++) {
    $ws2.Cells.Item(1, $i + 1).Value = $hotelHeaders[$i]
}

$ws2.Cells.Item(2, 1).Value = 32712
$ws2.Cells.Item(2, 2).Value = "Quality Inn New Orleans"
$ws2.Cells.Item(2, 3).Value = "Louisiana"
$ws2.Cells.Item(2, 4).Value = "New Orleans"
$ws2.Cells.Item(2, 5).Value = 70128
$ws2.Cells.Item(2, 6).Value = "https://www.tripadvisor.com/Hotel_Review-g60864-d93132-Reviews-Quality_Inn-New_Orleans_Louisiana.html"
$ws2.Cells.Item(2, 7).Value = "Quality Inn"

# English_Reviews_num, Local_Rank and Total_Reviews_num are stored as *text*
# (not numbers) in the source data, e.g. "9" and "139". Force a text format
# on those cells before assigning so they aren't auto-coerced to numbers,
# then restore the default "Normal" style so no stray styling is left behind.
$txtRange = $ws2.Range("H2:J2")
$txtRange.NumberFormat = "@"
$ws2.Cells.Item(2, 8).Value = "9"
$ws2.Cells.Item(2, 9).Value = "139"
$ws2.Cells.Item(2, 10).Value = "9"
$txtRange.Style = "Normal"
